$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the existing B26 entry (shared string index 51) in-place.
# Its text changes from "2019.5.28  13：44-" to "2019.5.29  18：00-21：00"
$ws.Range("B26").Value = "2019.5.29  18：00-21：00"
$ws.Range("C26").Value = "完成对学生表格的读取并存入数据库操作"

# New row 27 created below, pushing previous content down conceptually by
# inserting a fresh row before row 26 with the new log entry, then the
# original (now-edited) entry ends up on row 27.
$ws.Rows("26").Insert()

$ws.Range("B26").Value = "2019.5.28  13：44-17：30"
$ws.Range("C26").Value = "完成student实体的构建并学习poi的基本操作"

$ws.Range("C26").Select()
